# Fix TSP - out of errors.
# Update the "Fitness" column (C) values for rows 2-12 on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value  = 4235.486797775921
$ws.Range("C3").Value  = 4235.486797775921
$ws.Range("C4").Value  = 4235.486797775921
$ws.Range("C5").Value  = 4235.486797775921
$ws.Range("C6").Value  = 4162.080602175031
$ws.Range("C7").Value  = 4162.080602175031
$ws.Range("C8").Value  = 3927.023624666128
$ws.Range("C9").Value  = 3899.788819976822
$ws.Range("C10").Value = 3899.788819976822
$ws.Range("C11").Value = 3882.811179538926
$ws.Range("C12").Value = 3882.811179538926
